$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 15
$ws.Range("B1").Value = 3.072109460830688
$ws.Range("C1").Value = 2.773976564407349
$ws.Range("D1").Value = 3.071811199188232
$ws.Range("E1").Value = 15
